# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (D) and "Correspond Handback
# DateTime" (G) columns for the second data row (c27a0754-... entry) on
# each language sheet, reflecting a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-13 08:31:33"
$wsZhCn.Range("G3").Value = "2016-01-13 08:32:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-13 08:31:53"
$wsDeDe.Range("G3").Value = "2016-01-13 08:33:18"
